# finish 1st industry level code
# Add 3 new rows (12-14) to the "Changes from matlab" sheet, mirroring the
# existing file/function/line/changes layout (A=file, B=function, C=line, D=changes).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Changes from matlab")

# Row 12: D = "compress = true"
$ws.Range("D12").Value = "compress = true"

# Row 13: B = "add_time", A = Tools path (order matches shared-string insertion order)
$ws.Range("B13").Value = "add_time"
$ws.Range("A13").Value = "C:\Users\e0375379\Downloads\DT\Validus_SMECombined\Validus_SMECombined\ProdCode\Industry_Level\ZF_Code\Tools"

# Row 14: B = "get_individual_first_use_time", D = "julia trans date type, which is very slow"
$ws.Range("B14").Value = "get_individual_first_use_time"
$ws.Range("D14").Value = "julia trans date type, which is very slow"

# Update the selection/active cell to D14 on this sheet
$ws.Activate()
$ws.Range("D14").Select()

# Update the selection on the "Julia tip" sheet to C29 (view-only change)
$ws2 = $wb.Worksheets.Item("Julia tip")
$ws2.Activate()
$ws2.Range("C29").Select()

# Re-activate the first sheet so it is the one shown/tab-selected, matching the diff
$ws.Activate()
